$d = $word.ActiveDocument

# 1) Merge the split runs in the first paragraph: "Na min" + "h" + "a trajetória..." -> "Na minha trajetória..."
$d.Content.Find.Execute("Na minha trajetória", $true, $false, $false, $false, $false, $true, 1, $false, "Na minha trajetória", 2)

# 2) Merge the split runs "ABRASME \x{2013} " + "Associação Brasileira de Saúde Mental" into one run.
$d.Content.Find.Execute("ABRASME " + [char]8211 + " Associação Brasileira de Saúde Mental", $true, $false, $false, $false, $false, $true, 1, $false, "ABRASME " + [char]8211 + " Associação Brasileira de Saúde Mental", 2)

# 3) Add suppressAutoHyphens to Normal style paragraph format
$d.Styles("Normal").ParagraphFormat.Hyphenation = $false
